# Sync automático del tracker (cada 3h)
# Appends the newly finished/settled matches to the results tracker sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("14316315", "2025-08-12", "Reilly Opelka", "Francisco Comesaña", "Gana Francisco Comesaña", 2.75),
    @("14316409", "2025-08-12", "Ella Seidel", "McCartney Kessler", "Gana Ella Seidel", 3.75),
    @("14316466", "2025-08-12", "Iva Jovic", "Barbora Krejcikova", "Gana Iva Jovic", 2.1),
    @("14399578", "2025-08-12", "Hady Habib", "Alexander Shevchenko", "Gana Hady Habib", 2.63),
    @("14399582", "2025-08-12", "Mark Lajal", "Jaime Faria", "Gana Jaime Faria", 2.5),
    @("14399576", "2025-08-12", "Shintaro Mochizuki", "Alex Rybakov", "Gana Alex Rybakov", 4),
    @("14399572", "2025-08-12", "Nikoloz Basilashvili", "Christopher Eubanks", "Gana Christopher Eubanks", 2.2)
)

$startRow = 172

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    # event_id -> keep as plain text so it round-trips exactly like the
    # source feed (string, not a parsed number)
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $row[0]

    # fecha -> force text so Excel doesn't auto-convert the ISO date string
    # into a date serial number
    $ws.Cells.Item($r, 2).NumberFormat = "@"
    $ws.Cells.Item($r, 2).Value = $row[1]

    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]

    $ws.Cells.Item($r, 6).Value = $row[5]

    # resultado / profit -> still pending, leave as blank placeholders
    # (match hasn't been settled yet by the sync job)
    $ws.Cells.Item($r, 7).NumberFormat = "@"
    $ws.Cells.Item($r, 7).Value = ""
    $ws.Cells.Item($r, 8).NumberFormat = "@"
    $ws.Cells.Item($r, 8).Value = ""
}
